$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '96.028.38'
$ws.Range("E2").Value = '  -1.08%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.428.21'
$ws.Range("E3").Value = '  +2.99%  '

# Row 4
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.90'
$ws.Range("E5").Value = '  -1.38%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '643.71'
$ws.Range("E6").Value = '  -1.23%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.40'
$ws.Range("E7").Value = '  +2.90%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.405'
$ws.Range("E8").Value = '  -1.91%  '

# Row 9
$ws.Range("E9").Value = '  +0.10%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.989'
$ws.Range("E10").Value = '  +0.54%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.416.85'
$ws.Range("E11").Value = '  +2.73%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.29'
$ws.Range("E12").Value = '  +8.27%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.199'
$ws.Range("E13").Value = '  -2.39%  '

# Row 14
$ws.Range("B14").Value = 'WrappedBTC'
$ws.Range("C14").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '96.187.69'
$ws.Range("E14").Value = '  -0.65%  '

# Row 15
$ws.Range("B15").Value = 'Toncoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.06'
$ws.Range("E15").Value = '  +1.12%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.088.28'
$ws.Range("E16").Value = '  +3.50%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000250'
$ws.Range("E17").Value = '  -0.26%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.56'
$ws.Range("E18").Value = '  -0.69%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.443.53'
$ws.Range("E19").Value = '  +3.54%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.02'
$ws.Range("E20").Value = '  +7.30%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.64'
$ws.Range("E21").Value = '  +11.63%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.494'
$ws.Range("E22").Value = '  +2.54%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '510.05'
$ws.Range("E23").Value = '  +3.80%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.26'
$ws.Range("E24").Value = '  -0.42%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000193'
$ws.Range("E25").Value = '  -1.72%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.47'
$ws.Range("E26").Value = '  +1.01%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '92.09'
$ws.Range("E27").Value = '  -0.43%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.24'
$ws.Range("E28").Value = '  +2.18%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.633.45'
$ws.Range("E29").Value = '  +3.86%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.84'
$ws.Range("E30").Value = '  +9.91%  '

# Row 31
$ws.Range("E31").Value = '  +0.13%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.73'
$ws.Range("E32").Value = '  +11.06%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.138'
$ws.Range("E33").Value = '  -2.72%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.183'
$ws.Range("E34").Value = '  -1.43%  '

# Row 35
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.11%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '30.41'
$ws.Range("E36").Value = '  +8.79%  '

# Row 37
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.574'
$ws.Range("E37").Value = '  +5.45%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.73'
$ws.Range("E38").Value = '  +2.77%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.45'
$ws.Range("E39").Value = '  -1.51%  '

# Row 40
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  -0.14%  '

# Row 41
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.150'
$ws.Range("E41").Value = '  +0.41%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '502.48'
$ws.Range("E42").Value = '  -0.02%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.893'
$ws.Range("E43").Value = '  +8.41%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.23'
$ws.Range("E44").Value = '  -1.14%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.69'
$ws.Range("E45").Value = '  +3.53%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0415'
$ws.Range("E46").Value = '  +2.34%  '

# Row 47
$ws.Range("B47").Value = 'MantraDAO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.56'
$ws.Range("E47").Value = '  -3.02%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.47'
$ws.Range("E48").Value = '  +0.84%  '

# Row 49
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.29'
$ws.Range("E49").Value = '  +5.70%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.18'
$ws.Range("E50").Value = '  +10.91%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.30'
$ws.Range("E51").Value = '  -1.02%  '
